$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = '27.609.53'
$c.NumberFormat = "General"
$ws.Cells.Item(2, 5).Value = '  +0.07%  '

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = '1.844.28'
$c.NumberFormat = "General"
$ws.Cells.Item(3, 5).Value = '  +0.12%  '

$ws.Cells.Item(4, 5).Value = '  +0.24%  '

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '312.83'
$c.NumberFormat = "General"

$ws.Cells.Item(6, 5).Value = '  +0.23%  '

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = '0.4284'
$c.NumberFormat = "General"
$ws.Cells.Item(7, 5).Value = '  +1.01%  '

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = '0.3635'
$c.NumberFormat = "General"
$ws.Cells.Item(8, 5).Value = '  -0.14%  '

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '0.07316'
$c.NumberFormat = "General"
$ws.Cells.Item(9, 5).Value = '  +0.81%  '

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '0.8803'
$c.NumberFormat = "General"
$ws.Cells.Item(10, 5).Value = '  -1.44%  '

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '20.65'
$c.NumberFormat = "General"
$ws.Cells.Item(11, 5).Value = '  +0.07%  '

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '1.916.98'
$c.NumberFormat = "General"
$ws.Cells.Item(12, 5).Value = '  +4.61%  '

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '5.346'
$c.NumberFormat = "General"
$ws.Cells.Item(13, 5).Value = '  -0.17%  '

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '6.533'
$c.NumberFormat = "General"
$ws.Cells.Item(14, 5).Value = '  -0.58%  '

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = '0.06951'
$c.NumberFormat = "General"
$ws.Cells.Item(15, 5).Value = '  +1.10%  '

$ws.Cells.Item(16, 5).Value = '  +0.28%  '

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '79.57'
$c.NumberFormat = "General"
$ws.Cells.Item(17, 5).Value = '  +1.48%  '

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = '0.000008974'
$c.NumberFormat = "General"
$ws.Cells.Item(18, 5).Value = '  +1.58%  '

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.NumberFormat = "General"
$ws.Cells.Item(19, 5).Value = '  +0.15%  '

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '15.37'
$c.NumberFormat = "General"
$ws.Cells.Item(20, 5).Value = '  -0.64%  '

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '27.544.44'
$c.NumberFormat = "General"
$ws.Cells.Item(21, 5).Value = '  -0.13%  '

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '4.993'
$c.NumberFormat = "General"
$ws.Cells.Item(22, 5).Value = '  +0.18%  '

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '10.32'
$c.NumberFormat = "General"
$ws.Cells.Item(23, 5).Value = '  -2.34%  '

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '2.151.61'
$c.NumberFormat = "General"
$ws.Cells.Item(24, 5).Value = '  +4.39%  '

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '1.987'
$c.NumberFormat = "General"
$ws.Cells.Item(25, 5).Value = '  -1.68%  '

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '155.44'
$c.NumberFormat = "General"
$ws.Cells.Item(26, 5).Value = '  +0.36%  '

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '18.52'
$c.NumberFormat = "General"
$ws.Cells.Item(27, 5).Value = '  -0.29%  '

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '119.75'
$c.NumberFormat = "General"
$ws.Cells.Item(28, 5).Value = '  +1.08%  '

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = '5.218'
$c.NumberFormat = "General"
$ws.Cells.Item(29, 5).Value = '  -0.22%  '

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '1.887'
$c.NumberFormat = "General"
$ws.Cells.Item(30, 5).Value = '  +2.86%  '

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '0.08899'
$c.NumberFormat = "General"
$ws.Cells.Item(31, 5).Value = '  -0.09%  '

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '0.7656'
$c.NumberFormat = "General"
$ws.Cells.Item(32, 5).Value = '  -1.90%  '

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '2.966'
$c.NumberFormat = "General"
$ws.Cells.Item(33, 5).Value = '  +0.10%  '

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = '4.526'
$c.NumberFormat = "General"
$ws.Cells.Item(34, 5).Value = '  -0.95%  '

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = '1.131'
$c.NumberFormat = "General"
$ws.Cells.Item(35, 5).Value = '  +2.64%  '

$ws.Cells.Item(36, 5).Value = '  +0.18%  '

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = '0.05484'
$c.NumberFormat = "General"
$ws.Cells.Item(37, 5).Value = '  +1.45%  '

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = '1.105'
$c.NumberFormat = "General"
$ws.Cells.Item(38, 5).Value = '  +0.74%  '

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '0.01938'
$c.NumberFormat = "General"
$ws.Cells.Item(39, 5).Value = '  +0.85%  '

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '2.823'
$c.NumberFormat = "General"
$ws.Cells.Item(40, 5).Value = '  +1.55%  '

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '0.1667'
$c.NumberFormat = "General"
$ws.Cells.Item(41, 5).Value = '  +1.07%  '

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '0.5079'
$c.NumberFormat = "General"
$ws.Cells.Item(42, 5).Value = '  +0.25%  '

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '6.580'
$c.NumberFormat = "General"
$ws.Cells.Item(43, 5).Value = '  -3.87%  '

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '8.417'
$c.NumberFormat = "General"
$ws.Cells.Item(44, 5).Value = '  +2.36%  '

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '0.06548'
$c.NumberFormat = "General"
$ws.Cells.Item(45, 5).Value = '  -0.97%  '

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '105.89'
$c.NumberFormat = "General"
$ws.Cells.Item(46, 5).Value = '  +0.84%  '

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '10.35'
$c.NumberFormat = "General"
$ws.Cells.Item(47, 5).Value = '  +0.27%  '

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '0.4655'
$c.NumberFormat = "General"
$ws.Cells.Item(48, 5).Value = '  -1.11%  '

$ws.Cells.Item(49, 5).Value = '  +0.29%  '

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '1.640'
$c.NumberFormat = "General"
$ws.Cells.Item(50, 5).Value = '  +0.69%  '

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '64.75'
$c.NumberFormat = "General"
$ws.Cells.Item(51, 5).Value = '  +0.28%  '
